# Generate Report for Archive
# - Flip localization status from "Ready for handoff" to "In Translation"
#   on every sheet (Overview's per-language columns + each language
#   sheet's Status column).
# - The two affected columns on "Overview" and the Status column on each
#   language sheet then re-size (the shorter text needs a narrower
#   column) to match the refreshed report.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    [void]$ws.Cells.Replace("Ready for handoff", "In Translation")
}

$overview = $wb.Worksheets.Item("Overview")
$overview.Columns(5).ColumnWidth = 12.5
$overview.Columns(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns(3).ColumnWidth = 12.5
